$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = 0.023809523809523808
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0.09523809523809523
$ws.Range("I5").Value = 4
$ws.Range("L5").Value = 0.047619047619047616
$ws.Range("M5").Value = 2
$ws.Range("D6").Value = 0.037037037037037035
$ws.Range("E6").Value = 1
$ws.Range("H8").Value = 0.07547169811320754
$ws.Range("I8").Value = 4
$ws.Range("N8").Value = 0.05660377358490566
$ws.Range("O8").Value = 3
$ws.Range("L9").Value = 0.2
$ws.Range("M9").Value = 6
$ws.Range("H10").Value = 0.12820512820512819
$ws.Range("I10").Value = 5
$ws.Range("L10").Value = 0.05128205128205128
$ws.Range("M10").Value = 2
$ws.Range("J11").Value = 0.06666666666666667
$ws.Range("K11").Value = 1
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0.14285714285714285
$ws.Range("O13").Value = 3
$ws.Range("N14").Value = 0.13043478260869565
$ws.Range("O14").Value = 3
$ws.Range("C15").Value = 32
$ws.Range("D15").Value = 0.03125
$ws.Range("F15").Value = 0.0625
$ws.Range("H15").Value = 0.125
$ws.Range("L15").Value = 0.0625
$ws.Range("N15").Value = 0.09375
$ws.Range("F16").Value = 0.0967741935483871
$ws.Range("G16").Value = 3
$ws.Range("L19").Value = 0.26315789473684209
$ws.Range("M19").Value = 5
$ws.Range("D20").Value = 0.10714285714285714
$ws.Range("E20").Value = 3
$ws.Range("H20").Value = 0.21428571428571427
$ws.Range("I20").Value = 6
$ws.Range("J22").Value = 0.05263157894736842
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0.10526315789473684
$ws.Range("M22").Value = 2
$ws.Range("J29").Value = 0.06896551724137931
$ws.Range("K29").Value = 2
$ws.Range("D32").Value = 0.05357142857142857
$ws.Range("E32").Value = 3
$ws.Range("F41").Value = 0.14285714285714285
$ws.Range("G41").Value = 5
$ws.Range("L41").Value = 0.02857142857142857
$ws.Range("M41").Value = 1
$ws.Range("N41").Value = 0.08571428571428572
$ws.Range("O41").Value = 3
$ws.Range("L43").Value = 0.0625
$ws.Range("M43").Value = 2
$ws.Range("D45").Value = 0.08
$ws.Range("E45").Value = 2
$ws.Range("F45").Value = 0.12
$ws.Range("G45").Value = 3
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("N45").Value = 0.16
$ws.Range("O45").Value = 4
$ws.Range("D46").Value = 0.09523809523809523
$ws.Range("E46").Value = 2
$ws.Range("L48").Value = 0.1875
$ws.Range("M48").Value = 9
$ws.Range("D49").Value = 0.05555555555555555
$ws.Range("E49").Value = 3
$ws.Range("H49").Value = 0.24074074074074073
$ws.Range("I49").Value = 13
$ws.Range("J49").Value = 0.018518518518518517
$ws.Range("K49").Value = 1
$ws.Range("N49").Value = 0.05555555555555555
$ws.Range("O49").Value = 3
$ws.Range("D50").Value = 0.04878048780487805
$ws.Range("E50").Value = 2
$ws.Range("L50").Value = 0.0975609756097561
$ws.Range("M50").Value = 4
$ws.Range("N50").Value = 0.17073170731707318
$ws.Range("O50").Value = 7
$ws.Range("F51").Value = 0.12903225806451613
$ws.Range("G51").Value = 4
$ws.Range("L51").Value = 0.12903225806451613
$ws.Range("M51").Value = 4
$ws.Range("N52").Value = 0.11764705882352941
$ws.Range("O52").Value = 2
$ws.Range("H54").Value = 0.12903225806451613
$ws.Range("I54").Value = 4
$ws.Range("N54").Value = 0.19354838709677419
$ws.Range("O54").Value = 6
